$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
# D-column cells are forced to Text format first so numeric-looking strings
# (e.g. "191.59") are not auto-converted to actual numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.466.12"
$ws.Range("E2").Value = "  +2.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.357.65"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.59"
$ws.Range("E5").Value = "  +3.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "591.00"
$ws.Range("E6").Value = "  +1.67%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.76"
$ws.Range("E10").Value = "  +2.70%  "

$ws.Range("E11").Value = "  +1.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.946.98"
$ws.Range("E12").Value = "  +2.78%  "

$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.47"
$ws.Range("E14").Value = "  +3.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.510.36"
$ws.Range("E15").Value = "  +2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "453.92"
$ws.Range("E18").Value = "  +13.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.82"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.73"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.84"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.18"
$ws.Range("E22").Value = "  +6.61%  "

$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.519.49"
$ws.Range("E24").Value = "  +3.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.523"
$ws.Range("E25").Value = "  +1.97%  "

$ws.Range("E26").Value = "  +2.91%  "

$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.43"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.29"
$ws.Range("E31").Value = "  +2.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.55"
$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.95"
$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +5.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.10"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.20"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("E41").Value = "  +0.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.50"
$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.714.87"
$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "335.33"
$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.44"
$ws.Range("E50").Value = "  +5.16%  "

$ws.Range("E51").Value = "  +3.59%  "

# Rows 45 and 46 swap: InjectiveProtocol moves to row 45 (with updated price/volume),
# Hedera moves to row 46 (values unchanged)
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.40"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0688"
$ws.Range("E46").Value = "  +0.15%  "
